# 自动更新Excel文件
# Advances the tracking sheet by one day: each row's remaining-day counter
# (column E) ticks down by 1 relative to its total-day allotment (column D)
# and start date (column F). Rows whose counter has bottomed out at 1 (i.e.
# expiring the next day) are auto-renewed: the counter resets to the full
# allotment and the start date resets to the new "today".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newToday = 20260217

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $dCell = $ws.Cells.Item($r, 4)
    $eCell = $ws.Cells.Item($r, 5)
    $fCell = $ws.Cells.Item($r, 6)

    $dVal = $dCell.Value2
    $eVal = $eCell.Value2
    $fVal = $fCell.Value2

    if ($eVal -eq $null -or $dVal -eq $null -or $fVal -eq $null) {
        continue
    }

    # Skip rows whose start date isn't a well-formed YYYYMMDD value (data
    # glitch elsewhere in the sheet) - those are left untouched upstream too.
    $fStr = [string]$fVal
    if ($fStr.Length -ne 8) {
        continue
    }

    if ($eVal -le 1) {
        # Expiring tomorrow -> auto-renew starting today.
        $eCell.Value = $dVal
        $fCell.Value = $newToday
    } else {
        $eCell.Value = $eVal - 1
    }
}
